# Update the "Seasonality Index" (column P) values on the
# "Forecast Comparison" sheet to reflect the corrected week.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$newValues = @{
    2  = 0.99
    3  = 1.09
    4  = 0.95
    5  = 0.84
    6  = 1.04
    7  = 1.11
    8  = 1.13
    9  = 1.08
    10 = 1.06
    11 = 1.09
    12 = 1.08
    13 = 1.16
    14 = 0.9
    15 = 0.96
    16 = 0.97
    17 = 0.91
}

foreach ($row in $newValues.Keys) {
    $ws.Range("P$row").Value = $newValues[$row]
}
